# Add new "intervention_type" column (K) with per-row intervention type values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1 -- copy header style from J1 (bold, border, centered) then set text
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Cells.Item(1, 11).Value = "intervention_type"

# Row values for K2:K90, in order
$values = @(
    "PROCEDURE",
    "PROCEDURE",
    "DRUG",
    "PROCEDURE",
    "DRUG",
    "DRUG",
    "DEVICE",
    "OTHER",
    "PROCEDURE",
    "DEVICE",
    "GENETIC",
    "DRUG",
    "OTHER",
    "RADIATION",
    "BIOLOGICAL",
    "DEVICE",
    "PROCEDURE",
    "PROCEDURE",
    "DRUG",
    "DRUG",
    "PROCEDURE",
    "GENETIC",
    "DEVICE",
    "DRUG",
    "BEHAVIORAL",
    "DRUG",
    "DEVICE",
    "DRUG",
    "DEVICE",
    "OTHER",
    "PROCEDURE",
    "BEHAVIORAL",
    "DRUG",
    "DRUG",
    "BEHAVIORAL",
    "DRUG",
    "DRUG",
    "BIOLOGICAL",
    "PROCEDURE",
    "OTHER",
    "DRUG",
    "PROCEDURE",
    "OTHER",
    "DRUG",
    "DEVICE",
    "OTHER",
    "DEVICE",
    "PROCEDURE",
    "DRUG",
    "DRUG",
    "PROCEDURE",
    "OTHER",
    "DRUG",
    "DRUG",
    "GENETIC",
    "BEHAVIORAL",
    "DRUG",
    "BEHAVIORAL",
    "DRUG",
    "OTHER",
    "OTHER",
    "BEHAVIORAL",
    "OTHER",
    "OTHER",
    "DRUG",
    "DEVICE",
    "DRUG",
    "DRUG",
    "BEHAVIORAL",
    "DRUG",
    "DEVICE",
    "OTHER",
    "DEVICE",
    "DEVICE",
    "OTHER",
    "DEVICE",
    "BEHAVIORAL",
    "BIOLOGICAL",
    "BIOLOGICAL",
    "DRUG",
    "RADIATION",
    "DRUG",
    "DIETARY_SUPPLEMENT",
    "OTHER",
    "PROCEDURE",
    "BEHAVIORAL",
    "PROCEDURE",
    "OTHER",
    "DEVICE",
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}

Write-Output "Done. K1=$($ws.Cells.Item(1,11).Value()) K90=$($ws.Cells.Item(90,11).Value())"
